$d = $word.ActiveDocument

$replacements = @(
    @{old="622÷3="; new="350÷4="},
    @{old="773÷5="; new="453÷2="},
    @{old="464÷4="; new="418÷9="},
    @{old="757÷3="; new="848÷7="},
    @{old="202÷5="; new="535÷4="},
    @{old="842÷3="; new="993÷2="},
    @{old="180÷4="; new="663÷2="},
    @{old="125÷7="; new="183÷8="},
    @{old="687÷6="; new="732÷6="},
    @{old="202÷8="; new="760÷7="},
    @{old="423÷2="; new="762÷2="},
    @{old="142÷6="; new="176÷3="},
    @{old="344÷5="; new="903÷4="},
    @{old="816÷8="; new="944÷6="},
    @{old="215÷4="; new="371÷4="},
    @{old="955÷8="; new="136÷8="},
    @{old="792÷4="; new="809÷8="},
    @{old="937÷5="; new="761÷2="},
    @{old="727÷9="; new="107÷5="},
    @{old="626÷9="; new="264÷4="},
    @{old="771÷4="; new="769÷2="},
    @{old="666÷7="; new="244÷8="},
    @{old="150÷4="; new="686÷3="},
    @{old="272÷4="; new="836÷5="},
    @{old="507÷7="; new="488÷3="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
